# Daily refresh of COVID country stats ("Pais" sheet).
# Source data is kept sorted descending by column B (Casos totales),
# so a handful of neighbouring rows swap country labels where the
# day-over-day totals changed the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 14:47"

# Row 5: India
$ws.Range("B5").Value = 7602414
$ws.Range("C5").Value = 7678
$ws.Range("D5").Value = 6737145
$ws.Range("E5").Value = 749973
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 115296

# Row 18: Irak
$ws.Range("B18").Value = 434598
$ws.Range("C18").Value = 3920
$ws.Range("D18").Value = 366134
$ws.Range("E18").Value = 58098
$ws.Range("G18").Value = 49
$ws.Range("H18").Value = 10366

# Row 25: Arabia Saudita
$ws.Range("B25").Value = 342968
$ws.Range("C25").Value = 385
$ws.Range("D25").Value = 329270
$ws.Range("E25").Value = 8481
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 5217

# Row 29: Paises Bajos
$ws.Range("B29").Value = 244391
$ws.Range("C29").Value = 8165
$ws.Range("G29").Value = 46
$ws.Range("H29").Value = 6814

# Row 42: Kuwait
$ws.Range("A42").Value = "Kuwait"
$ws.Range("B42").Value = 117718
$ws.Range("C42").Value = 886
$ws.Range("D42").Value = 109198
$ws.Range("E42").Value = 7806
$ws.Range("H42").Value = 714

# Row 43: Emiratos Arabes Unidos
$ws.Range("A43").Value = "Emiratos Arabes Unidos"
$ws.Range("B43").Value = 117594
$ws.Range("C43").Value = 1077
$ws.Range("D43").Value = 110313
$ws.Range("E43").Value = 6811
$ws.Range("G43").Value = 4
$ws.Range("H43").Value = 470

# Row 46: Suecia
$ws.Range("A46").Value = "Suecia"
$ws.Range("B46").Value = 106380
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("G46").Value = 6
$ws.Range("H46").Value = 5922

# Row 47: Egipto
$ws.Range("A47").Value = "Egipto"
$ws.Range("B47").Value = 105547
$ws.Range("D47").Value = 98314
$ws.Range("E47").Value = 1103
$ws.Range("H47").Value = 6130

# Row 56: Suiza
$ws.Range("E56").Value = 29422
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 2145

# Row 72: Estado de Palestina
$ws.Range("B72").Value = 48129
$ws.Range("C72").Value = 513
$ws.Range("D72").Value = 41455
$ws.Range("E72").Value = 6253
$ws.Range("G72").Value = 8
$ws.Range("H72").Value = 421

# Row 80: Dinamarca
$ws.Range("A80").Value = "Dinamarca"
$ws.Range("B80").Value = 36373
$ws.Range("C80").Value = 529
$ws.Range("D80").Value = 29998
$ws.Range("E80").Value = 5687
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 688

# Row 81: Serbia
$ws.Range("A81").Value = "Serbia"
$ws.Range("B81").Value = 36282
$ws.Range("D81").Value = 31536
$ws.Range("E81").Value = 3968
$ws.Range("H81").Value = 778

# Row 82: Bosnia y Herzegovina
$ws.Range("B82").Value = 35389
$ws.Range("C82").Value = 728
$ws.Range("D82").Value = 25560
$ws.Range("E82").Value = 8812
$ws.Range("G82").Value = 20
$ws.Range("H82").Value = 1017

# Row 110: Uganda
$ws.Range("B110").Value = 10788
$ws.Range("C110").Value = 97
$ws.Range("D110").Value = 7066
$ws.Range("E110").Value = 3625

# Row 135: Guinea Ecuatorial
$ws.Range("B135").Value = 5074
$ws.Range("C135").Value = 4
$ws.Range("E135").Value = 37

# Row 142: Islandia
$ws.Range("A142").Value = "Islandia"
$ws.Range("B142").Value = 4193
$ws.Range("C142").Value = 92
$ws.Range("D142").Value = 2930
$ws.Range("E142").Value = 1252
$ws.Range("H142").Value = 11

# Row 143: Mayotte
$ws.Range("A143").Value = "Mayotte"
$ws.Range("B143").Value = 4159
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 2964
$ws.Range("E143").Value = 1152
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 43

# Row 144: Estonia
$ws.Range("A144").Value = "Estonia"
$ws.Range("B144").Value = 4127
$ws.Range("C144").Value = 42
$ws.Range("D144").Value = 3270
$ws.Range("E144").Value = 786
$ws.Range("G144").Value = 3
$ws.Range("H144").Value = 71

# Row 182: Islas Feroe
$ws.Range("B182").Value = 488
$ws.Range("C182").Value = 3
$ws.Range("E182").Value = 16

# Row 190: Liechtenstein
$ws.Range("A190").Value = "Liechtenstein"
$ws.Range("C190").Value = 11
$ws.Range("D190").Value = 142
$ws.Range("E190").Value = 92

# Row 191: Islas Caimanes
$ws.Range("A191").Value = "Islas Caimanes"
$ws.Range("B191").Value = 235
$ws.Range("D191").Value = 212
$ws.Range("E191").Value = 22
